# Update "想去人数" (people interested) counts in both the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1936
    $ws.Range("F4").Value = 1191
    $ws.Range("F5").Value = 1326
    $ws.Range("F7").Value = 6047
}
